{"js": "// Update the 20x5 grid of addition/subtraction problems in the single\n// table of the document. Each table cell holds exactly one run whose\n// text is a math expression like \"54-38=\". The new values below are the\n// post-edit text for each cell, in row-major (top-left to bottom-right)\n// order, taken from the target diff.\nconst newValues = [\n  \"50-18=\", \"47-18=\", \"41+23=\", \"48-46=\", \"53-24=\",\n  \"45+44=\", \"60-38=\", \"83-33=\", \"46+28=\", \"70+17=\",\n  \"90-4=\", \"98-20=\", \"74-67=\", \"27-1=\", \"41-23=\",\n  \"36+32=\", \"76-8=\", \"65-8=\", \"37-14=\", \"96-69=\",\n  \"59-42=\", \"38+21=\", \"6+71=\", \"70-10=\", \"96-16=\",\n  \"22-4=\", \"96-8=\", \"71-31=\", \"22-0=\", \"56-17=\",\n  \"71-8=\", \"84-69=\", \"34+11=\", \"31+26=\", \"17+28=\",\n  \"88-4=\", \"1+88=\", \"69-11=\", \"91-34=\", \"2+27=\",\n  \"41+54=\", \"16+25=\", \"31+12=\", \"66-23=\", \"23-22=\",\n  \"12+22=\", \"32-8=\", \"44+28=\", \"51-32=\", \"3+29=\",\n  \"65-25=\", \"43-3=\", \"73-32=\", \"75+14=\", \"17+47=\",\n  \"82-64=\", \"45-40=\", \"68-5=\", \"73+10=\", \"80+0=\",\n  \"87-35=\", \"36+44=\", \"62-37=\", \"59-34=\", \"10+56=\",\n  \"81-65=\", \"29-1=\", \"23+11=\", \"19+60=\", \"63+34=\",\n  \"89-29=\", \"75-31=\", \"31-22=\", \"18-6=\", \"80+12=\",\n  \"59+7=\", \"37+41=\", \"75-57=\", \"16+55=\", \"43+23=\",\n  \"39+6=\", \"78+20=\", \"82-20=\", \"39+35=\", \"6-2=\",\n  \"2+19=\", \"96-31=\", \"24+25=\", \"73-69=\", \"74+5=\",\n  \"26-13=\", \"76-19=\", \"63-36=\", \"73-62=\", \"22-19=\",\n  \"81+0=\", \"14+22=\", \"44-24=\", \"30+25=\", \"74-32=\",\n];\n\nconst ROWS = 20;\nconst COLS = 5;\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nfor (let r = 0; r < ROWS; r++) {\n  for (let c = 0; c < COLS; c++) {\n    const cell = table.getCell(r, c);\n    cell.value = newValues[r * COLS + c];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the 20x5 grid of addition/subtraction problems in the single\n# table of the document. Each table cell holds exactly one run whose\n# text is a math expression like \"54-38=\". The new values below are the\n# post-edit text for each cell, in row-major (top-left to bottom-right)\n# order, taken from the target diff. Word COM table cells are 1-based.\n$newValues = @(\n  \"50-18=\", \"47-18=\", \"41+23=\", \"48-46=\", \"53-24=\",\n  \"45+44=\", \"60-38=\", \"83-33=\", \"46+28=\", \"70+17=\",\n  \"90-4=\", \"98-20=\", \"74-67=\", \"27-1=\", \"41-23=\",\n  \"36+32=\", \"76-8=\", \"65-8=\", \"37-14=\", \"96-69=\",\n  \"59-42=\", \"38+21=\", \"6+71=\", \"70-10=\", \"96-16=\",\n  \"22-4=\", \"96-8=\", \"71-31=\", \"22-0=\", \"56-17=\",\n  \"71-8=\", \"84-69=\", \"34+11=\", \"31+26=\", \"17+28=\",\n  \"88-4=\", \"1+88=\", \"69-11=\", \"91-34=\", \"2+27=\",\n  \"41+54=\", \"16+25=\", \"31+12=\", \"66-23=\", \"23-22=\",\n  \"12+22=\", \"32-8=\", \"44+28=\", \"51-32=\", \"3+29=\",\n  \"65-25=\", \"43-3=\", \"73-32=\", \"75+14=\", \"17+47=\",\n  \"82-64=\", \"45-40=\", \"68-5=\", \"73+10=\", \"80+0=\",\n  \"87-35=\", \"36+44=\", \"62-37=\", \"59-34=\", \"10+56=\",\n  \"81-65=\", \"29-1=\", \"23+11=\", \"19+60=\", \"63+34=\",\n  \"89-29=\", \"75-31=\", \"31-22=\", \"18-6=\", \"80+12=\",\n  \"59+7=\", \"37+41=\", \"75-57=\", \"16+55=\", \"43+23=\",\n  \"39+6=\", \"78+20=\", \"82-20=\", \"39+35=\", \"6-2=\",\n  \"2+19=\", \"96-31=\", \"24+25=\", \"73-69=\", \"74+5=\",\n  \"26-13=\", \"76-19=\", \"63-36=\", \"73-62=\", \"22-19=\",\n  \"81+0=\", \"14+22=\", \"44-24=\", \"30+25=\", \"74-32=\"\n)\n\n$ROWS = 20\n$COLS = 5\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\nfor ($r = 1; $r -le $ROWS; $r++) {\n  for ($c = 1; $c -le $COLS; $c++) {\n    $idx = (($r - 1) * $COLS) + ($c - 1)\n    $cell = $table.Cell($r, $c)\n    $cell.Range.Text = $newValues[$idx]\n  }\n}\n"}
